$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "sexo" section-header row (row 5) and the "cor ou raça" section-header
# row (row 8, which becomes row 7 once row 5 is removed) are removed
# entirely. Deleting the rows shifts the data below them up, which
# reproduces the row renumbering seen in the diff (old rows 6,7,9,10 become
# new rows 5,6,7,8) and automatically drops the now-unused "sexo" and
# "cor ou raça" shared strings.
$ws.Rows("5").Delete()
$ws.Rows("7").Delete()
